$d = $word.ActiveDocument

# Locate the paragraph that holds the "video presentation link" bullet
# ("Ссылка на видеопрезентацию: …") and remove it entirely, including its
# paragraph mark, so the list collapses back to just the two remaining
# bullets (presentation link, GitHub repo link).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*видеопрезентац*") {
        $target = $p
    }
}

if ($target -ne $null) {
    $target.Range.Delete()
}
